# Generate Report for Handback
# - Updates the "In Translation" status to "Handed back: in sync with en-US"
#   everywhere it appears (Overview + per-locale sheets).
# - Records the handback xliff file name + handback datetime for each
#   locale/file row (zh-cn and de-de sheets), and adds a hyperlinked
#   "Latest Target File" entry matching the existing "File Name" link.
# - Column widths on the affected (now wider) columns are refreshed.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$urlMd1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d4790f3e3a62a08c0ea7a5dcd5a57f067e4e4252/e2e/6414ea35-e7e0-4b2e-a564-bc1f13bed67b.md"
$urlMd2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d4790f3e3a62a08c0ea7a5dcd5a57f067e4e4252/e2e/e3fab7e3-f3b9-4397-aee8-e1ff1033ec76.md"

$md1 = "6414ea35-e7e0-4b2e-a564-bc1f13bed67b.md"
$md2 = "e3fab7e3-f3b9-4397-aee8-e1ff1033ec76.md"

# ----------------------------------------------------------------------
# Overview sheet: refresh the displayed status text for both file rows.
# ----------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = 29.16666667
$overview.Columns.Item(6).ColumnWidth = 29.16666667

# ----------------------------------------------------------------------
# zh-cn sheet
# ----------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

# Latest Target File (column I) now links back to the source markdown file,
# same as column A.
$zhcn.Range("I2").Value = $md1
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $urlMd1, [Type]::Missing, [Type]::Missing, $md1) | Out-Null

$zhcn.Range("I3").Value = $md2
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $urlMd2, [Type]::Missing, [Type]::Missing, $md2) | Out-Null

# Latest Handback File (column J) now has the generated xliff file name.
$zhcn.Range("J2").Value = "6414ea35-e7e0-4b2e-a564-bc1f13bed67b.18480d8d6b162324ebaa9c2416f07ebbf0660f16.zh-cn.xlf"
$zhcn.Range("J3").Value = "e3fab7e3-f3b9-4397-aee8-e1ff1033ec76.b8ac4ac10105f2893c27bc3123c3e5e78525758b.zh-cn.xlf"

# Latest Handback DateTime (column K) — handback just completed.
$zhcn.Range("K2").Value = "2016-08-17 18:23:51"
$zhcn.Range("K3").Value = "2016-08-17 18:23:51"

$zhcn.Columns.Item(3).ColumnWidth = 29.16666667
$zhcn.Columns.Item(9).ColumnWidth = 39.16666667
$zhcn.Columns.Item(10).ColumnWidth = 39.16666667

# ----------------------------------------------------------------------
# de-de sheet
# ----------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

$dede.Range("I2").Value = $md1
$dede.Hyperlinks.Add($dede.Range("I2"), $urlMd1, [Type]::Missing, [Type]::Missing, $md1) | Out-Null

$dede.Range("I3").Value = $md2
$dede.Hyperlinks.Add($dede.Range("I3"), $urlMd2, [Type]::Missing, [Type]::Missing, $md2) | Out-Null

$dede.Range("J2").Value = "6414ea35-e7e0-4b2e-a564-bc1f13bed67b.18480d8d6b162324ebaa9c2416f07ebbf0660f16.de-de.xlf"
$dede.Range("J3").Value = "e3fab7e3-f3b9-4397-aee8-e1ff1033ec76.b8ac4ac10105f2893c27bc3123c3e5e78525758b.de-de.xlf"

$dede.Range("K2").Value = "2016-08-17 18:23:58"
$dede.Range("K3").Value = "2016-08-17 18:23:58"

$dede.Columns.Item(3).ColumnWidth = 29.16666667
$dede.Columns.Item(9).ColumnWidth = 39.16666667
$dede.Columns.Item(10).ColumnWidth = 39.16666667
